$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New price/volume data (open, close, high, low, shares_outstanding) per row 2..14
$data = @{
    2  = @(56.09999847412109, 53.52000045776367, 69, 49.84999847412109, 1059263999)
    3  = @(56.09999847412109, 53.52000045776367, 69, 49.84999847412109, 1059263999)
    4  = @(56.09999847412109, 53.52000045776367, 69, 49.84999847412109, 1059263999)
    5  = @(56.09999847412109, 53.52000045776367, 69, 49.84999847412109, 1059263999)
    6  = @(56.09999847412109, 53.52000045776367, 69, 49.84999847412109, 1059263999)
    7  = @(53.47499847412109, 49.29000091552734, 56.33000183105469, 46.5, 1059263999)
    8  = @(73.5, 70.66999816894531, 79.94000244140625, 65.11000061035156, 1059263999)
    9  = @(126.3600006103516, 101.2099990844727, 135.4100036621094, 85.61000061035156, 1059263999)
    10 = @(163.4299926757812, 144.1699981689453, 188.75, 129.1000061035156, 1059263999)
    11 = @(142.8399963378906, 141.3000030517578, 164.1600036621094, 132.8000030517578, 1059263999)
    12 = @(124.995002746582, 159.5500030517578, 182.8800048828125, 124.995002746582, 1059263999)
    13 = @(105.7200012207031, 114.0500030517578, 114.8300018310547, 80, 1059263999)
    14 = @(160, 141.375, 166.8800048828125, 139.1100006103516, 1059263999)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
    $ws.Cells.Item($row, 9).Value = "ARM"
}
